$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Hermes" (A2) was replaced by a new ticker "BIIB.O" which now sits in B2,
# and the "x" marker that used to be in D2 moved over to C2.
$ws.Range("A2").ClearContents()
$ws.Range("D2").ClearContents()

$ws.Range("B2").Value = "BIIB.O"
$ws.Range("C2").Value = "x"

# Update the active selection to match the new focus cell.
$ws.Range("C2").Select()

# Set page setup (paper size / orientation) - also opts the sheet into
# printing as A4 portrait, matching the saved pageSetup element.
$ps = $ws.PageSetup
$ps.PaperSize = 9
$ps.Orientation = 1
